$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A512").Value = "Buying Opportunity"
$ws.Range("B512").Value = "support Zone"
$ws.Range("C512").Value = "long buildup"
$ws.Range("D512").Value = "Short buildup"
$ws.Range("E512").Value = "FII ENTERING"

$ws.Range("A513").Value = "ASHOKAMET"
$ws.Range("B513").Value = "AHLADA"
$ws.Range("D513").Value = "ALKEM"
$ws.Range("F513").Value = 22.15
$ws.Range("G513").Value = 138.55
$ws.Range("I513").Value = 4946.8

$ws.Range("A514").Value = "BALPHARMA"
$ws.Range("B514").Value = "ALKALI"
$ws.Range("D514").Value = "BAJAJFINSV"
$ws.Range("F514").Value = 121.3
$ws.Range("G514").Value = 110.4
$ws.Range("I514").Value = 1524.1

$ws.Range("A515").Value = "ICEMAKE"
$ws.Range("B515").Value = "AMBER"
$ws.Range("D515").Value = "BAJFINANCE"
$ws.Range("F515").Value = 715.35
$ws.Range("G515").Value = 3510.25
$ws.Range("I515").Value = 6616.45

$ws.Range("B516").Value = "ANANTRAJ"
$ws.Range("D516").Value = "BALKRISIND"
$ws.Range("G516").Value = 370.45
$ws.Range("I516").Value = 3091.3

$ws.Range("B517").Value = "APOLLO"
$ws.Range("D517").Value = "INFY"
$ws.Range("G517").Value = 106.55
$ws.Range("I517").Value = 1427.45

$ws.Range("B518").Value = "ATFL"
$ws.Range("D518").Value = "LTIM"
$ws.Range("G518").Value = 673.85
$ws.Range("I518").Value = 4773.15

$ws.Range("B519").Value = "BALKRISHNA"
$ws.Range("D519").Value = "MANAPPURAM"
$ws.Range("G519").Value = 25.2
$ws.Range("I519").Value = 167.95

$ws.Range("B520").Value = "BALMLAWRIE"
$ws.Range("D520").Value = "NTPC"
$ws.Range("G520").Value = 268.25
$ws.Range("I520").Value = 359.7

$ws.Range("B521").Value = "BORORENEW"
$ws.Range("D521").Value = "RELIANCE"
$ws.Range("G521").Value = 473.45
$ws.Range("I521").Value = 2849.7

$ws.Range("B522").Value = "CAREERP"
$ws.Range("G522").Value = 375.05

$ws.Range("B523").Value = "CCL"
$ws.Range("G523").Value = 582.25

$ws.Range("B524").Value = "COASTCORP"
$ws.Range("G524").Value = 219.75

$ws.Range("B525").Value = "CYIENT"
$ws.Range("G525").Value = 1732.7

$ws.Range("B526").Value = "DCW"
$ws.Range("G526").Value = 50.4

$ws.Range("B527").Value = "DECCANCE"
$ws.Range("G527").Value = 568.15

$ws.Range("B528").Value = "DISHTV"
$ws.Range("G528").Value = 15.5

$ws.Range("B529").Value = "EKC"
$ws.Range("G529").Value = 125.65

$ws.Range("B530").Value = "EMUDHRA"
$ws.Range("G530").Value = 744.2

$ws.Range("B531").Value = "EVERESTIND"
$ws.Range("G531").Value = 1083.1

$ws.Range("B532").Value = "FACT"
$ws.Range("G532").Value = 683.15

$ws.Range("B533").Value = "FOSECOIND"
$ws.Range("G533").Value = 3548.85

$ws.Range("B534").Value = "GFLLIMITED"
$ws.Range("G534").Value = 70.55

$ws.Range("B535").Value = "GINNIFILA"
$ws.Range("G535").Value = 31.6

$ws.Range("B536").Value = "GOACARBON"
$ws.Range("G536").Value = 798.25

$ws.Range("B537").Value = "GODREJPROP"
$ws.Range("G537").Value = 2721.75

$ws.Range("B538").Value = "GOLDIAM"
$ws.Range("G538").Value = 161.4

$ws.Range("B539").Value = "GRAVITA"
$ws.Range("G539").Value = 1077.95

$ws.Range("B540").Value = "IMAGICAA"
$ws.Range("G540").Value = 78.1

$ws.Range("B541").Value = "JBMA"
$ws.Range("G541").Value = 2039.75

$ws.Range("B542").Value = "JINDRILL"
$ws.Range("G542").Value = 670.5

$ws.Range("B543").Value = "JMFINANCIL"
$ws.Range("G543").Value = 79.85

$ws.Range("B544").Value = "JSL"
$ws.Range("G544").Value = 707.65

$ws.Range("B545").Value = "KECL"
$ws.Range("G545").Value = 181.55

$ws.Range("B546").Value = "KOTARISUG"
$ws.Range("G546").Value = 54.05

$ws.Range("B547").Value = "LICI"
$ws.Range("G547").Value = 993.3

$ws.Range("B548").Value = "MAHASTEEL"
$ws.Range("G548").Value = 102.1

$ws.Range("B549").Value = "MICEL"
$ws.Range("G549").Value = 49

$ws.Range("B550").Value = "MOLDTECH"
$ws.Range("G550").Value = 234.2

$ws.Range("B551").Value = "NAVNETEDUL"
$ws.Range("G551").Value = 148.65

$ws.Range("B552").Value = "NELCO"
$ws.Range("G552").Value = 722.6

$ws.Range("B553").Value = "NEWGEN"
$ws.Range("G553").Value = 851.3

$ws.Range("B554").Value = "NLCINDIA"
$ws.Range("G554").Value = 218.2

$ws.Range("B555").Value = "NSIL"
$ws.Range("G555").Value = 3405.4

$ws.Range("B556").Value = "NYKAA"
$ws.Range("G556").Value = 161.5

$ws.Range("B557").Value = "OMINFRAL"
$ws.Range("G557").Value = 108.4

$ws.Range("B558").Value = "PATINTLOG"
$ws.Range("G558").Value = 20.7

$ws.Range("B559").Value = "PFS"
$ws.Range("G559").Value = 38.9

$ws.Range("B560").Value = "PGHL"
$ws.Range("G560").Value = 4773.9

$ws.Range("B561").Value = "PILITA"
$ws.Range("G561").Value = 12.65

$ws.Range("B562").Value = "PRICOLLTD"
$ws.Range("G562").Value = 432.7

$ws.Range("B563").Value = "RAILTEL"
$ws.Range("G563").Value = 406.8

$ws.Range("B564").Value = "SAKSOFT"
$ws.Range("G564").Value = 245.45

$ws.Range("A565").Value = "30/05/2024"
